# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 39 (shifting the existing rows 39-96
# down to 40-97) and populate the new row with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 39..96 down to 40..97, creating a new blank row 39.
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with the boilerplate (same on every
# row of this sheet) plus the new reading's Fecha / Volumen / Precio values.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44868
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100108
$ws.Range("H39").Value = "Tropicales y subtropicales"
$ws.Range("I39").Value = 100108007
$ws.Range("J39").Value = "Coco"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 15
$ws.Range("N39").Value = 32000
$ws.Range("O39").Value = 32000
$ws.Range("P39").Value = 32000
$ws.Range("Q39").Value = "$/malla 20 unidades"
$ws.Range("R39").Value = "Perú"
$ws.Range("S39").Value = 1600
$ws.Range("T39").Value = 20
